$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLink = "https://www.fiercebiotech.com/medtech/promega-msi-tech-wins-fda-approval-companion-diagnostic-keytruda-lenvima-combo-uterine"
$newKeyword = "companion diagnostic"
$newTitle = '<a href="https://www.fiercebiotech.com/medtech/promega-msi-tech-wins-fda-approval-companion-diagnostic-keytruda-lenvima-combo-uterine" hreflang="en">Promega wins FDA approval as companion diagnostic for Keytruda-Lenvima combo in uterine cancer</a>'

$rowNum = 63
$linkCell = $ws.Cells.Item($rowNum, 1)
$keywordCell = $ws.Cells.Item($rowNum, 2)
$titleCell = $ws.Cells.Item($rowNum, 3)

# Populate the new row's values first.
$linkCell.Value = $newLink
$keywordCell.Value = $newKeyword
$titleCell.Value = $newTitle

# Turn the link cell into a live hyperlink (mirrors every other row in column A).
$ws.Hyperlinks.Add($linkCell, $newLink)

# Hyperlinks.Add() stamps its own style variant - reapply the same
# "Hyperlink" look the rest of column A uses so the new cell matches.
$ws.Range("A62").Copy()
$linkCell.PasteSpecial(-4122)
$excel.CutCopyMode = 0
